$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:C2").Value = 0.88
$ws.Range("D2").Value = 0.18

$ws.Range("A3:C3").Value = 0.97
$ws.Range("D3").Value = 0.79

$ws.Range("A4:C4").Value = 0.0000000000000384
$ws.Range("D4").Value = 0.00044
